$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Tags" column (E) with a Text number format, matching the
# order the unique tag values need to land in the shared string table.
$ws.Range("E1").Value = "Tags"
$ws.Range("E20").Value = "@Debug"
$ws.Range("E2").Value = "@allPages, @BecomePartnerPage"
$ws.Range("E3").Value = "@allPages, @DownloadsPage"
$ws.Range("E6").Value = "@allPages, @ContactUsPage"
$ws.Range("E4").Value = "@allPages, @DownloadsGetStartedPage"
$ws.Range("E5").Value = "@allPages, @DownloadsStartTodayPage"
$ws.Range("E16").Value = "@allPages, @CouchbaseVsMongoDbPage"
$ws.Range("E17").Value = "@allPages, @PricingFormPage"
$ws.Range("E18").Value = "@allPages, @ProfessionalServicesPage"
$ws.Range("E19").Value = "@allPages, @RegisterDealPage"
$ws.Range("E21").Value = "@allPages, @ProductPage"

$ws.Range("E7").Value = "@allPages, @DownloadsPage"
$ws.Range("E8").Value = "@allPages, @DownloadsPage"
$ws.Range("E9").Value = "@allPages, @DownloadsPage"
$ws.Range("E10").Value = "@allPages, @DownloadsGetStartedPage"
$ws.Range("E11").Value = "@allPages, @DownloadsGetStartedPage"
$ws.Range("E12").Value = "@allPages, @DownloadsGetStartedPage"
$ws.Range("E13").Value = "@allPages, @DownloadsStartTodayPage"
$ws.Range("E14").Value = "@allPages, @DownloadsStartTodayPage"
$ws.Range("E15").Value = "@allPages, @DownloadsStartTodayPage"
$ws.Range("E22").Value = "@allPages, @ProductPage"

# Formatting: header cell bold (matches the other header cells in row 1)
# and text number format so tag strings such as "@Debug" are stored as text.
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1:E22").NumberFormat = "@"

# Column width for the new Tags column.
$ws.Columns.Item(5).ColumnWidth = 39.5

# Restore the selection to where the user left off after entering the data.
[void]$ws.Range("E23").Select()
